$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2253539253539254
$ws.Range("H2").Value = 0.2253539253539254
$ws.Range("I2").Value = 0.2194337194337195
$ws.Range("J2").Value = 0.1853240040476682
$ws.Range("K2").Value = 14.17
$ws.Range("L2").Value = 0.1823680823680824
$ws.Range("M2").Value = 11.43
$ws.Range("N2").Value = 0.09930495221546481
$ws.Range("O2").Value = 0.8066337332392378
$ws.Range("P2").Value = 11.43
$ws.Range("Q2").Value = 0.09930495221546481
$ws.Range("R2").Value = 0.8066337332392378
$ws.Range("U2").Value = 13.23
$ws.Range("V2").Value = 0.1149435273675065
$ws.Range("W2").Value = 0.09649126942964464
$ws.Range("X2").Value = 0.06399564953766468
$ws.Range("Y2").Value = 0.03249561989197997
$ws.Range("Z2").Value = 0.9939874632211845
$ws.Range("AA2").Value = 0.1736170394444319
$ws.Range("AB2").Value = 0.06355563128679435
$ws.Range("AC2").Value = 0.1100614081576375
$ws.Range("AD2").Value = 3.97
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 3.97
$ws.Range("AG2").Value = -9.26
$ws.Range("AH2").Value = 0.03334173175443017
$ws.Range("AI2").Value = 0.0274987878368082
$ws.Range("AJ2").Value = -0.08749055177626607
$ws.Range("AK2").Value = -0.07061156016470946
$ws.Range("AL2").Value = 0.208
$ws.Range("AM2").Value = 0.208
$ws.Range("AN2").Value = 0.2195796460176991
$ws.Range("AO2").Value = 81.97115384615385
$ws.Range("AP2").Value = -0.5121681415929203
$ws.Range("AQ2").Value = 81.97115384615385

# Row 3
$ws.Range("G3").Value = 0.2397868561278863
$ws.Range("H3").Value = 0.2397868561278863
$ws.Range("I3").Value = 0.2362344582593251
$ws.Range("J3").Value = 0.1992664705164536
$ws.Range("K3").Value = 11
$ws.Range("L3").Value = 0.1953818827708703
$ws.Range("M3").Value = 9.35
$ws.Range("N3").Value = 0.09999999999999999
$ws.Range("O3").Value = 0.85
$ws.Range("P3").Value = 9.35
$ws.Range("Q3").Value = 0.09999999999999999
$ws.Range("R3").Value = 0.85
$ws.Range("U3").Value = 5.55
$ws.Range("V3").Value = 0.05935828877005347
$ws.Range("W3").Value = 0.1097804391217565
$ws.Range("X3").Value = 0.06496516210034954
$ws.Range("Y3").Value = 0.04481527702140695
$ws.Range("Z3").Value = 0.9847822284414902
$ws.Range("AA3").Value = 0.1962340788888637
$ws.Range("AB3").Value = 0.06408512559860889
$ws.Range("AC3").Value = 0.1321489532902548
$ws.Range("AD3").Value = 3.97
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 3.97
$ws.Range("AG3").Value = -1.58
$ws.Range("AH3").Value = 0.04073048117369447
$ws.Range("AI3").Value = 0.03753427247801835
$ws.Range("AJ3").Value = -0.01718885987815491
$ws.Range("AK3").Value = -0.01576531630413091
$ws.Range("AL3").Value = 0.208
$ws.Range("AM3").Value = 0.208
$ws.Range("AN3").Value = 0.2835714285714286
$ws.Range("AO3").Value = 63.9423076923077
$ws.Range("AP3").Value = -0.1128571428571428
$ws.Range("AQ3").Value = 63.9423076923077

# Row 4
$ws.Range("G4").Value = 0.1873831775700935
$ws.Range("H4").Value = 0.1873831775700935
$ws.Range("I4").Value = 0.1752336448598131
$ws.Range("J4").Value = 0.148177570093458
$ws.Range("K4").Value = 3.17
$ws.Range("L4").Value = 0.1481308411214953
$ws.Range("M4").Value = 2.08
$ws.Range("N4").Value = 0.0962962962962963
$ws.Range("O4").Value = 0.6561514195583596
$ws.Range("P4").Value = 2.08
$ws.Range("Q4").Value = 0.0962962962962963
$ws.Range("R4").Value = 0.6561514195583596
$ws.Range("U4").Value = 7.68
$ws.Range("V4").Value = 0.3555555555555555
$ws.Range("W4").Value = 0.0832020997375328
$ws.Range("X4").Value = 0.06302613697497982
$ws.Range("Y4").Value = 0.02017596276255298
$ws.Range("Z4").Value = 1.019047619047619
$ws.Range("AA4").Value = 0.151
$ws.Range("AB4").Value = 0.06302613697497982
$ws.Range("AC4").Value = 0.08797386302502021
$ws.Range("AG4").Value = -7.68
$ws.Range("AJ4").Value = -0.5517241379310344
$ws.Range("AK4").Value = -0.2483829236739974
$ws.Range("AP4").Value = -1.882352941176471
